$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Remove the duplicated "Contact" row (old row 11) so everything below shifts up
$ws.Rows(11).Delete()

# Version bump 5.0.0 -> 6.0.0
$ws.Range("B3").Value = "6.0.0"

# Date bump
$ws.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher now has a value
$ws.Range("B9").Value = "Alvearie Team"

# The remaining "Contact" row becomes "Jurisdiction" / "United States of America"
$ws.Range("A10").Value = "Jurisdiction"
$ws.Range("B10").Value = "United States of America"

# Case Sensitive now has a value of "true"
$ws.Range("B14").Value = "true"
